# ICDC_ProgramsPage.xlsx edit
#
# The "stdCntVal" text for the Comparative Oncology Program (cell D2 on the
# "programs" sheet) had a stray leading line break baked into it. Fix the
# text, then restate the sheet's formatting the way Excel does when you
# select everything and set the number format to Text: every used cell
# picks up numFmtId 49 ("@"), the yellow highlight on the index header
# (E1) and the word-wrap on the description column (D2:D4) are kept, and
# the two long description rows grow a little to fit.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("programs")

# --- Fix the COP standard content text: drop the leading blank line ---
$ws.Range("D2").Value2 = "The COP is a core resource for CCR investigators who are interested in the use of comparative cancer models."

# --- Keep the description column's word-wrap ---
$ws.Range("D2:D4").WrapText = $true

# --- Keep the yellow highlight on the stdCntIndex header cell ---
$ws.Range("E1").Interior.ColorIndex = 6

# --- Row heights: row 2 no longer needs the extra line, rows 3 & 4 grow a bit ---
$ws.Rows.Item(2).EntireRow.AutoFit()
$ws.Rows.Item(3).RowHeight = 45
$ws.Rows.Item(4).RowHeight = 45

# --- Select the whole sheet (Ctrl+A) and format it as Text, same as the
#     author did to normalize the numeric-looking index/count columns ---
$ws.Cells.Select()
$ws.Cells.NumberFormat = "@"
